# Thailand Premier League - add new match result (2024-04-19) and shift
# the five already-entered 2024-04-20/21 matches down by one row.
#
# Row 192 becomes the brand new match (Port FC vs Chiangrai Utd, 19/04/2024).
# The data that used to live in rows 192-196 moves down into rows 193-197
# (the "id" column A is positional - A<row> == <row>-2 - and therefore is
# left untouched everywhere).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-MatchRow($ws, $r, $data) {
    foreach ($col in $data.Keys) {
        $addr = "$col$r"
        $ws.Range($addr).Value = $data[$col]
    }
}

# Row 197 is brand new - it does not exist yet, so first clone the number
# formatting (style) of row 196's A/E cells (bold+border style, and the
# custom date/time style) before filling in values, so no new style entries
# get created in styles.xml.
$ws.Range("A196").Copy()
$ws.Range("A197").PasteSpecial(-4122)
$ws.Range("E196").Copy()
$ws.Range("E197").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 192: brand new match data -----------------------------------
$row192 = @{
    B  = 6992709
    C  = "Thailand Premier League"
    D  = "Thailand Premier League"
    E  = 45401.375
    F  = "Port FC"
    G  = "Chiangrai Utd"
    H  = 1
    I  = 2
    J  = "A"
    K  = 1.363
    L  = 4.5
    M  = 7.5
    N  = 1.4
    O  = 4.5
    P  = 6.5
    Q  = -1.25
    R  = 1.8
    S  = 2
    T  = 3
    U  = 1.9
    V  = 1.9
    W  = -1
    X  = -1
    Y  = 5.5
    Z  = -1
    AA = 1
    AB = 0
    AC = -0
}
Set-MatchRow $ws 192 $row192

# --- Row 193: (was row 192) -------------------------------------------
$row193 = @{
    B = 6992714
    C = "Thailand Premier League"
    D = "Thailand Premier League"
    E = 45402.33333333334
    F = "BG Pathum United"
    G = "Chonburi"
    K = 1.444
    L = 4.333
    M = 6
    N = 1.45
    O = 4.5
    P = 5.5
    Q = -1.25
    R = 1.95
    S = 1.85
    T = 3
    U = 1.95
    V = 1.85
    W = 0
    X = 0
    Y = 0
    Z = 0
    AA = 0
}
Set-MatchRow $ws 193 $row193

# --- Row 194: (was row 193) -------------------------------------------
$row194 = @{
    B = 6995900
    C = "Thailand Premier League"
    D = "Thailand Premier League"
    E = 45402.35416666666
    F = "Police Tero FC"
    G = "Uthai Thani FC"
    K = 3.3
    L = 3.6
    M = 1.95
    N = 3.1
    O = 3.6
    P = 2.05
    Q = 0.25
    R = 1.975
    S = 1.825
    T = 3
    U = 1.95
    V = 1.85
    W = 0
    X = 0
    Y = 0
    Z = 0
    AA = 0
}
Set-MatchRow $ws 194 $row194

# --- Row 195: (was row 194) -------------------------------------------
$row195 = @{
    B = 6992713
    C = "Thailand Premier League"
    D = "Thailand Premier League"
    E = 45402.375
    F = "Khonkaen United"
    G = "Trat FC"
    K = 2.1
    L = 3.75
    M = 2.875
    N = 2.05
    O = 3.8
    P = 2.9
    Q = -0.25
    R = 1.85
    S = 1.95
    T = 3
    U = 2
    V = 1.8
    W = 0
    X = 0
    Y = 0
    Z = 0
    AA = 0
}
Set-MatchRow $ws 195 $row195

# --- Row 196: (was row 195) -------------------------------------------
$row196 = @{
    B = 6992710
    C = "Thailand Premier League"
    D = "Thailand Premier League"
    E = 45402.41666666666
    F = "Ratchaburi FC"
    G = "Buriram United"
    K = 5.25
    L = 3.75
    M = 1.571
    N = 5
    O = 3.75
    P = 1.6
    Q = 1
    R = 1.8
    S = 2
    T = 2.75
    U = 1.95
    V = 1.85
    W = 0
    X = 0
    Y = 0
    Z = 0
    AA = 0
}
Set-MatchRow $ws 196 $row196

# --- Row 197: new row (was row 196) -----------------------------------
$row197 = @{
    A = 195
    B = 6992715
    C = "Thailand Premier League"
    D = "Thailand Premier League"
    E = 45403.375
    F = "Nakhon Pathom FC"
    G = "Bangkok United"
    K = 4.75
    L = 3.75
    M = 1.615
    N = 6.5
    O = 4
    P = 1.45
    Q = 1
    R = 2.025
    S = 1.775
    T = 2.75
    U = 1.9
    V = 1.9
    W = 0
    X = 0
    Y = 0
    Z = 0
    AA = 0
}
Set-MatchRow $ws 197 $row197
